# Update "想去人数" (number of people wanting to go) figures that changed
# between data scrapes, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 363
    $ws.Range("F7").Value = 241
    $ws.Range("F10").Value = 412
}
